# Javascript_Gradebook_Storyboard.pptx - "Grades page being worked on"
#
# Changes applied to slide 4 (the "Grades" page):
#   1. Nudge the "Teacher Course History" textbox a bit to the right.
#   2. Add a 6x5 "Student List" grades table.
#   3. Add a new "Student List" header textbox above the table (styled
#      like the other nav textboxes on the page).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# ---------------------------------------------------------------------
# Helper: PowerPoint's Shape.Left/Top/Width/Height are single-precision
# (float32) points, and converting back to EMU on save truncates rather
# than rounds. 1 pt = 12700 EMU, so naively sending emu/12700.0 can miss
# the intended EMU value by a hair. Probe with the real shape itself and
# nudge up in tiny sub-point steps until the round-tripped value lands
# exactly on the requested EMU amount.
# ---------------------------------------------------------------------
function EmuToPt {
    param(
        [double]$Emu,
        $ProbeShape,
        [string]$Axis = "Left"
    )
    $base = $Emu / 12700.0
    for ($k = 0; $k -lt 400; $k++) {
        $cand = $base + ($k * (1.0 / 12700.0) / 8.0)
        switch ($Axis) {
            "Left"   { $ProbeShape.Left = $cand }
            "Top"    { $ProbeShape.Top = $cand }
            "Width"  { $ProbeShape.Width = $cand }
            "Height" { $ProbeShape.Height = $cand }
        }
        $readBack = switch ($Axis) {
            "Left"   { $ProbeShape.Left }
            "Top"    { $ProbeShape.Top }
            "Width"  { $ProbeShape.Width }
            "Height" { $ProbeShape.Height }
        }
        $emuBack = [math]::Round($readBack * 12700.0)
        if ($emuBack -eq $Emu) {
            return $cand
        }
    }
    return $base
}

function SetShapeRectEmu {
    param($Shape, $Left, $Top, $Width, $Height)
    if ($null -ne $Left)   { $Shape.Left   = EmuToPt $Left   $Shape "Left" }
    if ($null -ne $Top)    { $Shape.Top    = EmuToPt $Top    $Shape "Top" }
    if ($null -ne $Width)  { $Shape.Width  = EmuToPt $Width  $Shape "Width" }
    if ($null -ne $Height) { $Shape.Height = EmuToPt $Height $Shape "Height" }
}

# ---------------------------------------------------------------------
# 1. Move "Teacher Course History" textbox (shape id 11 / positional 10)
#    221673 -> 251130 EMU on X; Y/size unchanged.
# ---------------------------------------------------------------------
$courseHistory = $s.Shapes.Item(10)
SetShapeRectEmu $courseHistory 251130 $null $null $null

# ---------------------------------------------------------------------
# 2. Add the grades table (6 rows x 5 columns) - becomes shape id 12,
#    name "Table 11".
# ---------------------------------------------------------------------
$tblShape = $s.Shapes.AddTable(6, 5, 1, 1, 1, 1)
SetShapeRectEmu $tblShape 2514602 1414780 6477000 5227320

$tbl = $tblShape.Table

$colWidths = @(1216338, 1216338, 1377322, 1055354, 1611648)
for ($c = 1; $c -le 5; $c++) {
    $col = $tbl.Columns.Item($c)
    $ptw = $colWidths[$c - 1] / 12700.0
    for ($k = 0; $k -lt 400; $k++) {
        $cand = $ptw + ($k * (1.0 / 12700.0) / 8.0)
        $col.Width = $cand
        $back = [math]::Round($col.Width * 12700.0)
        if ($back -eq $colWidths[$c - 1]) { break }
    }
}

for ($r = 1; $r -le 6; $r++) {
    $row = $tbl.Rows.Item($r)
    $pth = 871220 / 12700.0
    for ($k = 0; $k -lt 400; $k++) {
        $cand = $pth + ($k * (1.0 / 12700.0) / 8.0)
        $row.Height = $cand
        $back = [math]::Round($row.Height * 12700.0)
        if ($back -eq 871220) { break }
    }
}

# Header row text.
$headers = @("Student Name", "Email:", "Assignments", "Grade", "Comments")
for ($c = 1; $c -le 5; $c++) {
    $cell = $tbl.Cell(1, $c)
    $cell.Shape.TextFrame.TextRange.Text = $headers[$c - 1]
}
# "Assignments" header renders smaller (empty data cells below it keep
# the table's default run size - there's no text run to attach a size
# to until something is typed into them).
$tbl.Cell(1, 3).Shape.TextFrame.TextRange.Font.Size = 16

# ---------------------------------------------------------------------
# 3. Add the "Student List" textbox (shape id 13, name "TextBox 12").
#    Duplicate the "Teacher Info" textbox (positional 9) so the new
#    shape inherits the same accent4 style + 24pt run formatting, then
#    reposition / resize / retext it.
# ---------------------------------------------------------------------
$styleSource = $s.Shapes.Item(9)
$dupRange = $styleSource.Duplicate()
$studentList = $dupRange.Item(1)
$studentList.Name = "TextBox 12"
$studentList.TextFrame.TextRange.Text = "Student List"
SetShapeRectEmu $studentList 251130 4495800 1901952 461665
